# This script re-shuffles the data rows (2-15) of the sheet so that each
# row ends up carrying the D, K, L, M, N, O, P, Q, R, S, T values that used
# to belong to a different row, per the target permutation below.
# Columns A, B, C, E, F, G, H, I, J are identical across all rows and are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values (D, K, L, M, N, O, P, Q, R, S, T) for every
# data row (2..15) before any of them gets overwritten.
$before = @{}
for ($r = 2; $r -le 15; $r++) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# Target row -> source row mapping (which row's old data should land here).
$mapping = @{
    2  = 14
    3  = 15
    4  = 6
    5  = 3
    6  = 12
    7  = 13
    8  = 4
    9  = 5
    10 = 7
    11 = 8
    12 = 11
    13 = 2
    14 = 9
    15 = 10
}

foreach ($targetRow in ($mapping.Keys | Sort-Object)) {
    $sourceRow = $mapping[$targetRow]
    $vals = $before[$sourceRow]

    $ws.Cells.Item($targetRow, 4).Value  = $vals.D
    $ws.Cells.Item($targetRow, 11).Value = $vals.K
    $ws.Cells.Item($targetRow, 12).Value = $vals.L
    $ws.Cells.Item($targetRow, 13).Value = $vals.M
    $ws.Cells.Item($targetRow, 14).Value = $vals.N
    $ws.Cells.Item($targetRow, 15).Value = $vals.O
    $ws.Cells.Item($targetRow, 16).Value = $vals.P
    $ws.Cells.Item($targetRow, 17).Value = $vals.Q
    $ws.Cells.Item($targetRow, 18).Value = $vals.R
    $ws.Cells.Item($targetRow, 19).Value = $vals.S
    $ws.Cells.Item($targetRow, 20).Value = $vals.T
}
